$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Contenu du stage" breakdown (rows 16-23 on the pie-chart source range
# D16:D23 / E16:E23): the counts for C#, COBOL, ASSEMBLEUR and ANDROID were
# wrong (COBOL/ASSEMBLEUR/ANDROID all showed 0 while C# wrongly held the
# full total of 32). Fix the counts and their matching percentages.

# C# : 32 -> 0  (0 %)
$ws.Range("E16").Value = 0
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "0 %"

# COBOL : 0 -> 29  (90.63 %)
$ws.Range("E17").Value = 29
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "90.63 %"

# ASSEMBLEUR : 0 -> 2  (6.25 %)
$ws.Range("E19").Value = 2
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "6.25 %"

# ANDROID : 0 -> 1  (3.13 %)
$ws.Range("E20").Value = 1
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "3.13 %"
